# Auto-generated edit script applying the diff's cell value changes
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 275.75
$ws.Range("J29").Value = 500
$ws.Range("L29").Value = 1500
$ws.Range("N29").Value = -2062
$ws.Range("H38").Value = 915
$ws.Range("I38").Value = 246.66667
$ws.Range("J38").Value = 2920
$ws.Range("K38").Value = 740.00001
$ws.Range("L38").Value = 8760
$ws.Range("M38").Value = -368.00001
$ws.Range("N38").Value = -9504
$ws.Range("H58").Value = 21743.922
$ws.Range("I58").Value = 328.57144
$ws.Range("J58").Value = 25150.908
$ws.Range("K58").Value = 985.71432
$ws.Range("L58").Value = 75452.724
$ws.Range("M58").Value = -835.71432
$ws.Range("N58").Value = -75752.724
$ws.Range("H87").Value = 22894.625
$ws.Range("J87").Value = 22894.625
$ws.Range("L87").Value = 22894.625
$ws.Range("N87").Value = -25390.625
$ws.Range("H90").Value = 22894.625
$ws.Range("J90").Value = 22894.625
$ws.Range("L90").Value = 68683.875
$ws.Range("N90").Value = -81163.875
$ws.Range("H96").Value = 988.4286
$ws.Range("I96").Value = 945
$ws.Range("J96").Value = 1005.8
$ws.Range("K96").Value = 2835
$ws.Range("L96").Value = 3017.4
$ws.Range("M96").Value = -1462
$ws.Range("N96").Value = -5763.4
$ws.Range("H132").Value = 4002430.8
$ws.Range("I132").Value = 4257532
$ws.Range("J132").Value = 5845.3335
$ws.Range("K132").Value = 12772596
$ws.Range("L132").Value = 17536.0005
$ws.Range("M132").Value = -12770066
$ws.Range("N132").Value = -22596.0005
$ws.Range("H135").Value = 762.4286
$ws.Range("I135").Value = 557.5263
$ws.Range("J135").Value = 2709
$ws.Range("K135").Value = 5017.736699999999
$ws.Range("L135").Value = 24381
$ws.Range("M135").Value = -2482.736699999999
$ws.Range("N135").Value = -29451
$ws.Range("H138").Value = 4485.0273
$ws.Range("I138").Value = 1407.6522
$ws.Range("J138").Value = 5900.62
$ws.Range("K138").Value = 4222.9566
$ws.Range("L138").Value = 17701.86
$ws.Range("M138").Value = 917.0434000000005
$ws.Range("N138").Value = -27981.86

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1587.9
$ws.Range("I32").Value = 1334.4
$ws.Range("J32").Value = 3869.4
$ws.Range("K32").Value = 1334.4
$ws.Range("L32").Value = 3869.4
$ws.Range("M32").Value = -1047.4
$ws.Range("N32").Value = -4443.4
$ws.Range("H102").Value = 3683.3333
$ws.Range("I102").Value = 2575
$ws.Range("J102").Value = 5900
$ws.Range("K102").Value = 2575
$ws.Range("L102").Value = 5900
$ws.Range("M102").Value = -953
$ws.Range("N102").Value = -9144

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H105").Value = 1850.6666
$ws.Range("I105").Value = 1688.3334
$ws.Range("J105").Value = 2013
$ws.Range("K105").Value = 1688.3334
$ws.Range("L105").Value = 2013
$ws.Range("M105").Value = 58.66660000000002
$ws.Range("N105").Value = -5507
$ws.Range("H141").Value = 35030.668
$ws.Range("J141").Value = 29587.273
$ws.Range("L141").Value = 29587.273
$ws.Range("N141").Value = -39947.273

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 15997
$ws.Range("J133").Value = 15997
$ws.Range("L133").Value = 15997
$ws.Range("N133").Value = -21057
$ws.Range("H141").Value = 29420
$ws.Range("J141").Value = 29420
$ws.Range("L141").Value = 29420
$ws.Range("N141").Value = -39780

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 877.6857
$ws.Range("I5").Value = 619.96295
$ws.Range("J5").Value = 1747.5
$ws.Range("K5").Value = 1859.88885
$ws.Range("L5").Value = 5242.5
$ws.Range("M5").Value = -1747.88885
$ws.Range("N5").Value = -5466.5
$ws.Range("H34").Value = 5855.8335
$ws.Range("I34").Value = 196
$ws.Range("J34").Value = 7345.263
$ws.Range("K34").Value = 588
$ws.Range("L34").Value = 22035.789
$ws.Range("M34").Value = -504
$ws.Range("N34").Value = -22203.789
$ws.Range("H39").Value = 1625.7894
$ws.Range("I39").Value = 750
$ws.Range("J39").Value = 1728.8235
$ws.Range("K39").Value = 2250
$ws.Range("L39").Value = 5186.470499999999
$ws.Range("M39").Value = -1956
$ws.Range("N39").Value = -5774.470499999999
$ws.Range("H55").Value = 2860.7693
$ws.Range("J55").Value = 3294.5454
$ws.Range("L55").Value = 9883.636200000001
$ws.Range("N55").Value = -10237.6362
$ws.Range("H87").Value = 14425
$ws.Range("H90").Value = 14425
$ws.Range("H120").Value = 16612
$ws.Range("H131").Value = 1798.8462
$ws.Range("I131").Value = 2807.7778
$ws.Range("J131").Value = 1264.7059
$ws.Range("K131").Value = 8423.3334
$ws.Range("L131").Value = 3794.1177
$ws.Range("M131").Value = -3383.3334
$ws.Range("N131").Value = -13874.1177
$ws.Range("H132").Value = 4583.3335
$ws.Range("I132").Value = 2750
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 24750
$ws.Range("L132").Value = 49500
$ws.Range("M132").Value = -22220
$ws.Range("N132").Value = -54560
$ws.Range("H133").Value = 3417.7827
$ws.Range("J133").Value = 2604.2144
$ws.Range("L133").Value = 7812.6432
$ws.Range("N133").Value = -17932.6432
$ws.Range("H134").Value = 2932.182
$ws.Range("I134").Value = 1488.1428
$ws.Range("J134").Value = 3606.0667
$ws.Range("K134").Value = 4464.428400000001
$ws.Range("L134").Value = 10818.2001
$ws.Range("M134").Value = 605.5715999999993
$ws.Range("N134").Value = -20958.2001
$ws.Range("H135").Value = 877.6857
$ws.Range("I135").Value = 619.96295
$ws.Range("J135").Value = 1747.5
$ws.Range("K135").Value = 5579.66655
$ws.Range("L135").Value = 15727.5
$ws.Range("M135").Value = -3044.66655
$ws.Range("N135").Value = -20797.5
$ws.Range("H137").Value = 2026.5676
$ws.Range("I137").Value = 1376
$ws.Range("J137").Value = 2791.9412
$ws.Range("K137").Value = 4128
$ws.Range("L137").Value = 8375.8236
$ws.Range("M137").Value = 972
$ws.Range("N137").Value = -18575.8236
$ws.Range("H138").Value = 3282.1
$ws.Range("I138").Value = 1384.2
$ws.Range("J138").Value = 5180
$ws.Range("K138").Value = 4152.6
$ws.Range("L138").Value = 15540
$ws.Range("M138").Value = 987.3999999999996
$ws.Range("N138").Value = -25820
$ws.Range("H139").Value = 7148271.5
$ws.Range("I139").Value = 10002381
$ws.Range("J139").Value = 12996.9
$ws.Range("K139").Value = 30007143
$ws.Range("L139").Value = 38990.7
$ws.Range("M139").Value = -30002003
$ws.Range("N139").Value = -49270.7
$ws.Range("H140").Value = 5211419
$ws.Range("I140").Value = 13889672
$ws.Range("J140").Value = 4467.5
$ws.Range("K140").Value = 41669016
$ws.Range("L140").Value = 13402.5
$ws.Range("M140").Value = -41663836
$ws.Range("N140").Value = -23762.5
$ws.Range("H141").Value = 2944.4443
$ws.Range("I141").Value = 2357.1428
$ws.Range("K141").Value = 7071.428400000001
$ws.Range("M141").Value = -1891.428400000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 29525
$ws.Range("J137").Value = 29525
$ws.Range("L137").Value = 29525
$ws.Range("N137").Value = -39725
$ws.Range("H141").Value = 32375
$ws.Range("J141").Value = 32375
$ws.Range("L141").Value = 32375
$ws.Range("N141").Value = -42735

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992
$ws.Range("H122").Value = 296175.78
$ws.Range("I122").Value = 386338.22
$ws.Range("J122").Value = 3147.875
$ws.Range("K122").Value = 1159014.66
$ws.Range("L122").Value = 9443.625
$ws.Range("M122").Value = -1156564.66
$ws.Range("N122").Value = -14343.625
$ws.Range("H135").Value = 39349.75
$ws.Range("J135").Value = 39349.75
$ws.Range("L135").Value = 39349.75
$ws.Range("N135").Value = -49489.75
